# chore: data driven added to tc3, tc4, tc5
# Populates PartnerInfo (tc4) and OnlineDegree (tc5) sheets with data-driven rows.

$wb = $excel.ActiveWorkbook

# Unicode right single quotation mark used in one of the long descriptions below.
$rsquo = [char]0x2019

# --- PartnerInfo sheet: Partner Link | Logo Displayed | Partner Name -------
$wsPartner = $wb.Worksheets.Item("PartnerInfo")

$wsPartner.Range("A2").Value = "https://www.coursera.org/partners/huji"
$wsPartner.Range("B2").Value = "'true"
$wsPartner.Range("B2").Style = "Normal"
$wsPartner.Range("C2").Value = "Hebrew University of Jerusalem"

$wsPartner.Range("A3").Value = "https://www.coursera.org/partners/technion"
$wsPartner.Range("B3").Value = "'true"
$wsPartner.Range("B3").Style = "Normal"
$wsPartner.Range("C3").Value = "Technion - Israel Institute of Technology"

$wsPartner.Range("A4").Value = "https://www.coursera.org/partners/telaviv"
$wsPartner.Range("B4").Value = "'true"
$wsPartner.Range("B4").Style = "Normal"
$wsPartner.Range("C4").Value = "Tel Aviv University"

$wsPartner.Range("A5").Value = "https://www.coursera.org/partners/yadvashem"
$wsPartner.Range("B5").Value = "'true"
$wsPartner.Range("B5").Style = "Normal"
$wsPartner.Range("C5").Value = "Yad Vashem"

# --- OnlineDegree sheet: Card Details ---------------------------------------
$wsDegree = $wb.Worksheets.Item("OnlineDegree")

$wsDegree.Range("A2").Value = "Indian Institute of Technology Guwahati`nBachelor of Science in Data Science & AI`nNamed as one of the world${rsquo}s top universities for the study of Data Science (QS World University Rankings by Subject 2024)`nApplication due July 30, 2025"

$wsDegree.Range("A3").Value = "Birla Institute of Technology & Science, Pilani`nBachelor of Science in Computer Science`nRanked #7 among Technical Universities in India (The Week-Hansa Research Best Universities Survey 2024)`nApplication due August 3, 2025"

$wsDegree.Range("A4").Value = "University of London`nBachelor of Science in Computer Science`nSpecialise in ML and AI, data science, web and mobile development, physical computing and IoT, game development, VR, or UX`nApplication due September 8, 2025"

$wsDegree.Range("A5").Value = "University of London`nInternational Foundation Programme for Computer Science`nSecure your future in computer science regardless of your academic or professional experience and qualifications`nApplication due December 5, 2025"

# Multi-line values trigger Excel's auto row-height; AutoFit restores the
# default (non-custom) row height so the rows stay free of ht/customHeight.
$wsDegree.Rows(2).AutoFit()
$wsDegree.Rows(3).AutoFit()
$wsDegree.Rows(4).AutoFit()
$wsDegree.Rows(5).AutoFit()
